$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.6690309026450687
$ws.Range("D2").Value = 0.5104361041434156

# Row 3
$ws.Range("C3").Value = 1.967811484812883
$ws.Range("D3").Value = 0.06182177307170478
$ws.Range("G3").Value = "No"

# Row 4
$ws.Range("C4").Value = 0.5782936390784065
$ws.Range("D4").Value = 0.5689390361786382

# Row 5
$ws.Range("C5").Value = 2.095827268067037
$ws.Range("D5").Value = 0.04782190839436717

# Row 6
$ws.Range("C6").Value = 1.048891195603637
$ws.Range("D6").Value = 0.3056149475113965

# Row 7
$ws.Range("C7").Value = -0.1787734428045338
$ws.Range("D7").Value = 0.8597517230004268

# Row 8
$ws.Range("C8").Value = 1.467934355533703
$ws.Range("D8").Value = 0.1562747859273512

# Row 9
$ws.Range("C9").Value = -1.480458328694501
$ws.Range("D9").Value = 0.1529308425492046

# Row 10
$ws.Range("C10").Value = 0.1818847760516814
$ws.Range("D10").Value = 0.8573386186330474

# Row 11
$ws.Range("C11").Value = 1.379923119172511
$ws.Range("D11").Value = 0.181474198113639
